$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2
$ws.Range("A4").Value = 21.91
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = Get-Date -Year 2025 -Month 9 -Day 23
$ws.Range("H2").Value = "0x3c499c542cef5e3811e1192ce70d8cc03d5c3359"
